$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 3.655376239133008
$ws.Range("D2").Value = 5.56231141504255
$ws.Range("E2").Value = 16.40242787886047
$ws.Range("F2").Value = 32.36197923181384
$ws.Range("G2").Value = 3.645616317449174
$ws.Range("I2").Value = 25.99990240752329
$ws.Range("K2").Value = 17.96025305051723
$ws.Range("N2").Value = 17.47684377557205
$ws.Range("C3").Value = 3.663888575478786
$ws.Range("D3").Value = 5.599016170178225
$ws.Range("E3").Value = 15.47254267126582
$ws.Range("F3").Value = 32.03323250404111
$ws.Range("G3").Value = 3.649933239856683
$ws.Range("I3").Value = 25.8231708744707
$ws.Range("K3").Value = 17.35964762128669
$ws.Range("N3").Value = 17.55955343139331
$ws.Range("C4").Value = 3.669290399719855
$ws.Range("D4").Value = 5.622601975847249
$ws.Range("E4").Value = 14.87820719737089
$ws.Range("F4").Value = 31.84258871122776
$ws.Range("G4").Value = 3.65271492413717
$ws.Range("I4").Value = 25.72372004389993
$ws.Range("K4").Value = 16.9855477180453
$ws.Range("N4").Value = 17.61233771403611
$ws.Range("C5").Value = 3.671536245968598
$ws.Range("D5").Value = 5.632477152759714
$ws.Range("E5").Value = 14.63040422795807
$ws.Range("F5").Value = 31.7677907102282
$ws.Range("G5").Value = 3.653881593335942
$ws.Range("I5").Value = 25.6854973699724
$ws.Range("K5").Value = 16.83201961310351
$ws.Range("N5").Value = 17.63435313504398
$ws.Range("C6").Value = 3.671911873051576
$ws.Range("D6").Value = 5.634132850358406
$ws.Range("E6").Value = 14.58892715342619
$ws.Range("F6").Value = 31.75554713437465
$ws.Range("G6").Value = 3.654077322038594
$ws.Range("I6").Value = 25.67929034714063
$ws.Range("K6").Value = 16.8064691682354
$ws.Range("N6").Value = 17.6380393745426
$ws.Range("C7").Value = 3.669320506899929
$ws.Range("D7").Value = 5.622734087892296
$ws.Range("E7").Value = 14.87488754017481
$ws.Range("F7").Value = 31.84156816372516
$ws.Range("G7").Value = 3.652730523983326
$ws.Range("I7").Value = 25.72319519933962
$ws.Range("K7").Value = 16.98348119703811
$ws.Range("N7").Value = 17.61263257213921
$ws.Range("C8").Value = 3.658275225309121
$ws.Range("D8").Value = 5.574749501829905
$ws.Range("E8").Value = 16.08680457875361
$ws.Range("F8").Value = 32.24634689852449
$ws.Range("G8").Value = 3.647077677994069
$ws.Range("I8").Value = 25.93710101778329
$ws.Range("K8").Value = 17.75441949308822
$ws.Range("N8").Value = 17.50494839058367
$ws.Range("C9").Value = 3.63798345095792
$ws.Range("D9").Value = 5.488975153392017
$ws.Range("E9").Value = 18.32748483418707
$ws.Range("F9").Value = 33.12551447449925
$ws.Range("G9").Value = 3.637025632574907
$ws.Range("I9").Value = 26.42714049814342
$ws.Range("K9").Value = 19.21356627799341
$ws.Range("N9").Value = 17.3095377565088
$ws.Range("C10").Value = 3.623878636381131
$ws.Range("D10").Value = 5.431032797347356
$ws.Range("E10").Value = 19.97626867643695
$ws.Range("F10").Value = 33.81834443618711
$ws.Range("G10").Value = 3.630260601879611
$ws.Range("I10").Value = 26.82809771281628
$ws.Range("K10").Value = 20.24101895558892
$ws.Range("N10").Value = 17.17541399219309
$ws.Range("C11").Value = 3.617630143407416
$ws.Range("D11").Value = 5.405777147758234
$ws.Range("E11").Value = 20.68516370601517
$ws.Range("F11").Value = 34.14252003503439
$ws.Range("G11").Value = 3.627315611095931
$ws.Range("I11").Value = 27.01886109253611
$ws.Range("K11").Value = 20.69652537214765
$ws.Range("N11").Value = 17.11641328404041
$ws.Range("C12").Value = 3.615287656151758
$ws.Range("D12").Value = 5.396372361408026
$ws.Range("E12").Value = 20.94772850635074
$ws.Range("F12").Value = 34.2664665877734
$ws.Range("G12").Value = 3.626219305021642
$ws.Range("I12").Value = 27.09225167364421
$ws.Range("K12").Value = 20.86714528811873
$ws.Range("N12").Value = 17.0943580152426
$ws.Range("C13").Value = 3.615791106596433
$ws.Range("D13").Value = 5.398390772145644
$ws.Range("E13").Value = 20.89144131290244
$ws.Range("F13").Value = 34.23972134512016
$ws.Range("G13").Value = 3.626454576029558
$ws.Range("I13").Value = 27.07639526691429
$ws.Range("K13").Value = 20.830484796774
$ws.Range("N13").Value = 17.09909528547757
$ws.Range("C14").Value = 3.617436953541873
$ws.Range("D14").Value = 5.40500022143312
$ws.Range("E14").Value = 20.7068826626015
$ws.Range("F14").Value = 34.15269399392898
$ws.Range("G14").Value = 3.62722503944578
$ws.Range("I14").Value = 27.02487617788888
$ws.Range("K14").Value = 20.71060072001043
$ws.Range("N14").Value = 17.11459304546648
$ws.Range("C15").Value = 3.618448152957226
$ws.Range("D15").Value = 5.409069417733359
$ws.Range("E15").Value = 20.5930707571334
$ws.Range("F15").Value = 34.099538786371
$ws.Range("G15").Value = 3.627699427189351
$ws.Range("I15").Value = 26.99346781150977
$ws.Range("K15").Value = 20.63692020267204
$ws.Range("N15").Value = 17.12412317703654
$ws.Range("C16").Value = 3.624290328577983
$ws.Range("D16").Value = 5.432705535558168
$ws.Range("E16").Value = 19.92911604777793
$ws.Range("F16").Value = 33.79733092600375
$ws.Range("G16").Value = 3.630455716573326
$ws.Range("I16").Value = 26.81579497998433
$ws.Range("K16").Value = 20.21099714354884
$ws.Range("N16").Value = 17.17931012378741
$ws.Range("C17").Value = 3.627916978581825
$ws.Range("D17").Value = 5.447488214108918
$ws.Range("E17").Value = 19.51129036086107
$ws.Range("F17").Value = 33.61416724030541
$ws.Range("G17").Value = 3.63218042999954
$ws.Range("I17").Value = 26.70890624971491
$ws.Range("K17").Value = 19.94654200236462
$ws.Range("N17").Value = 17.21367934551228
$ws.Range("C18").Value = 3.630018757986602
$ws.Range("D18").Value = 5.456094573659968
$ws.Range("E18").Value = 19.2670930338906
$ws.Range("F18").Value = 33.50967034971094
$ws.Range("G18").Value = 3.633184915818202
$ws.Range("I18").Value = 26.64821733468431
$ws.Range("K18").Value = 19.79332395438793
$ws.Range("N18").Value = 17.23363721895799
$ws.Range("C19").Value = 3.630733116595556
$ws.Range("D19").Value = 5.459026347171392
$ws.Range("E19").Value = 19.18374530614908
$ws.Range("F19").Value = 33.47443938413442
$ws.Range("G19").Value = 3.633527164799156
$ws.Range("I19").Value = 26.62780639413427
$ws.Range("K19").Value = 19.74126145212948
$ws.Range("N19").Value = 17.24042725081518
$ws.Range("C20").Value = 3.627529280881727
$ws.Range("D20").Value = 5.445903831779709
$ws.Range("E20").Value = 19.55616939324701
$ws.Range("F20").Value = 33.63357770441638
$ws.Range("G20").Value = 3.63199554114328
$ws.Range("I20").Value = 26.72020326361186
$ws.Range("K20").Value = 19.97480985676576
$ws.Range("N20").Value = 17.21000107724877
$ws.Range("C21").Value = 3.61695288927401
$ws.Range("D21").Value = 5.403054547525046
$ws.Range("E21").Value = 20.76125124485752
$ws.Range("F21").Value = 34.17822464445155
$ws.Range("G21").Value = 3.626998223965317
$ws.Range("I21").Value = 27.03997771643891
$ws.Range("K21").Value = 20.74586553651288
$ws.Range("N21").Value = 17.11003320746459
$ws.Range("C22").Value = 3.610178434733113
$ws.Range("D22").Value = 5.375976972005397
$ws.Range("E22").Value = 21.51459597477855
$ws.Range("F22").Value = 34.54105865263053
$ws.Range("G22").Value = 3.623842271099748
$ws.Range("I22").Value = 27.25566198003501
$ws.Range("K22").Value = 21.23883043158344
$ws.Range("N22").Value = 17.04637025812803
$ws.Range("C23").Value = 3.613781625410042
$ws.Range("D23").Value = 5.390343811275343
$ws.Range("E23").Value = 21.1156426129839
$ws.Range("F23").Value = 34.34681379340025
$ws.Range("G23").Value = 3.625516639236346
$ws.Range("I23").Value = 27.13995214400629
$ws.Range("K23").Value = 20.97677757664773
$ws.Range("N23").Value = 17.08019620409248
$ws.Range("C24").Value = 3.627704506742857
$ws.Range("D24").Value = 5.44661979590998
$ws.Range("E24").Value = 19.53589199416734
$ws.Range("F24").Value = 33.62479971393388
$ws.Range("G24").Value = 3.632079089125183
$ws.Range("I24").Value = 26.71509350376529
$ws.Range("K24").Value = 19.96203362051194
$ws.Range("N24").Value = 17.21166340370804
$ws.Range("C25").Value = 3.643329641852024
$ws.Range("D25").Value = 5.511288347934626
$ws.Range("E25").Value = 17.70018660641737
$ws.Range("F25").Value = 32.87904558090899
$ws.Range("G25").Value = 3.639635365399462
$ws.Range("I25").Value = 26.28721686724907
$ws.Range("K25").Value = 18.82583464165293
$ws.Range("N25").Value = 17.36073124975275
